$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows: capacity values (column B)
$ws.Range("B2").Value = 20000
$ws.Range("B3").Value = 15000

# Add new rows 4-6, following the same pattern as rows 2-3
# Columns: A=ID_Battery, B=capacity, C=capacity_unit, D=charge_efficiency,
#          E=charge_power_max, F=charge_power_max_unit, G=discharge_efficiency,
#          H=discharge_power_max, I=discharge_power_max_unit

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 10000
$ws.Range("C4").Value = "Wh"
$ws.Range("D4").Value = 0.95
$ws.Range("E4").Value = 4500
$ws.Range("F4").Value = "W"
$ws.Range("G4").Value = 0.95
$ws.Range("H4").Value = 4500
$ws.Range("I4").Value = "W"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 5000
$ws.Range("C5").Value = "Wh"
$ws.Range("D5").Value = 0.95
$ws.Range("E5").Value = 4500
$ws.Range("F5").Value = "W"
$ws.Range("G5").Value = 0.95
$ws.Range("H5").Value = 4500
$ws.Range("I5").Value = "W"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = "Wh"
$ws.Range("D6").Value = 0.95
$ws.Range("E6").Value = 4500
$ws.Range("F6").Value = "W"
$ws.Range("G6").Value = 0.95
$ws.Range("H6").Value = 4500
$ws.Range("I6").Value = "W"

# Update selection & zoom to match the recorded view state
$ws.Range("B11").Select()
$excel.ActiveWindow.Zoom = 182
